$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number + date range) ---
$ws.Range("A8").Value2 = "Volume 30   Number  30"
$ws.Range("C9").Value2 = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# --- Cells that change type/style (copy format+value from a donor cell, then set final value) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("I14").Copy($ws.Range("D30"))
$ws.Range("D30").Value2 = 1
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("E30").Value2 = -100
$ws.Range("I14").Copy($ws.Range("G30"))
$ws.Range("G30").Value2 = 1
$ws.Range("K14").Copy($ws.Range("H30"))
$ws.Range("H30").Value2 = 0

# --- Plain value updates (style unchanged) ---
$ws.Range("L15").Value2 = 100
$ws.Range("M15").Value2 = 233.333333333333
$ws.Range("N15").Value2 = -28.571428571428
$ws.Range("D16").Value2 = 3
$ws.Range("E16").Value2 = -100
$ws.Range("F16").Value2 = 6
$ws.Range("G16").Value2 = 6
$ws.Range("H16").Value2 = 0
$ws.Range("J16").Value2 = 48
$ws.Range("K16").Value2 = -4.166666666666
$ws.Range("L16").Value2 = 17.948717948717
$ws.Range("M16").Value2 = -51.578947368421
$ws.Range("N16").Value2 = -89.176470588235
$ws.Range("G17").Value2 = 18
$ws.Range("H17").Value2 = -22.222222222222
$ws.Range("I17").Value2 = 121
$ws.Range("J17").Value2 = 123
$ws.Range("K17").Value2 = -1.626016260162
$ws.Range("L17").Value2 = 53.164556962025
$ws.Range("M17").Value2 = 40.697674418604
$ws.Range("N17").Value2 = -35.978835978836
$ws.Range("C18").Value2 = 4
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 9
$ws.Range("G18").Value2 = 20
$ws.Range("H18").Value2 = -55
$ws.Range("I18").Value2 = 69
$ws.Range("J18").Value2 = 104
$ws.Range("K18").Value2 = -33.653846153846
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = -68.202764976958
$ws.Range("N18").Value2 = -92.864529472595
$ws.Range("C19").Value2 = 10
$ws.Range("D19").Value2 = 19
$ws.Range("E19").Value2 = -47.368421052631
$ws.Range("F19").Value2 = 36
$ws.Range("G19").Value2 = 52
$ws.Range("H19").Value2 = -30.769230769230
$ws.Range("I19").Value2 = 312
$ws.Range("J19").Value2 = 349
$ws.Range("K19").Value2 = -10.601719197707
$ws.Range("L19").Value2 = 26.829268292682
$ws.Range("M19").Value2 = 36.842105263157
$ws.Range("N19").Value2 = -16.353887399463
$ws.Range("C20").Value2 = 2
$ws.Range("D20").Value2 = 5
$ws.Range("E20").Value2 = -60
$ws.Range("F20").Value2 = 10
$ws.Range("H20").Value2 = 42.857142857142
$ws.Range("I20").Value2 = 87
$ws.Range("J20").Value2 = 55
$ws.Range("K20").Value2 = 58.181818181818
$ws.Range("L20").Value2 = 85.106382978723
$ws.Range("M20").Value2 = 2.352941176470
$ws.Range("N20").Value2 = -91.951896392229
$ws.Range("C21").Value2 = 18
$ws.Range("D21").Value2 = 35
$ws.Range("E21").Value2 = -48.571428571428
$ws.Range("F21").Value2 = 75
$ws.Range("G21").Value2 = 104
$ws.Range("H21").Value2 = -27.884615384615
$ws.Range("I21").Value2 = 647
$ws.Range("J21").Value2 = 692
$ws.Range("K21").Value2 = -6.502890173410
$ws.Range("L21").Value2 = 33.127572016460
$ws.Range("M21").Value2 = -9.383753501400
$ws.Range("N21").Value2 = -78.821603927986
$ws.Range("L22").Value2 = -28.571428571428
$ws.Range("C24").Value2 = 19
$ws.Range("D24").Value2 = 23
$ws.Range("E24").Value2 = -17.391304347826
$ws.Range("F24").Value2 = 72
$ws.Range("G24").Value2 = 88
$ws.Range("H24").Value2 = -18.181818181818
$ws.Range("I24").Value2 = 647
$ws.Range("J24").Value2 = 643
$ws.Range("K24").Value2 = 0.622083981337
$ws.Range("L24").Value2 = 45.720720720720
$ws.Range("M24").Value2 = 26.862745098039
$ws.Range("C25").Value2 = 9
$ws.Range("D25").Value2 = 14
$ws.Range("E25").Value2 = -35.714285714285
$ws.Range("F25").Value2 = 31
$ws.Range("G25").Value2 = 31
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 204
$ws.Range("J25").Value2 = 191
$ws.Range("K25").Value2 = 6.806282722513
$ws.Range("L25").Value2 = 25.925925925925
$ws.Range("M25").Value2 = -21.839080459770
$ws.Range("L26").Value2 = 22.222222222222
$ws.Range("F27").Value2 = 4
$ws.Range("G27").Value2 = 9
$ws.Range("H27").Value2 = -55.555555555555
$ws.Range("L27").Value2 = -5.128205128205
$ws.Range("F30").Value2 = 1
$ws.Range("J30").Value2 = 10
$ws.Range("K30").Value2 = -50

Write-Host "Edit complete"